{"js": "// The Notation images were moved into their own \"Notation\" sub-folder on\n// the site, so every \"img/NotationN.png\" markdown-style image reference in\n// the body needs to become \"img/Notation/NotationN.png\".\nconst names = [\"Notation1.png\", \"Notation2.png\", \"Notation3.png\", \"Notation4.png\"];\n\nfor (const name of names) {\n  const oldRef = \"img/\" + name;\n  const newRef = \"img/Notation/\" + name;\n\n  const results = context.document.body.search(oldRef, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newRef, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The Notation images were moved into their own \"Notation\" sub-folder on\n# the site, so every \"img/NotationN.png\" markdown-style image reference in\n# the body needs to become \"img/Notation/NotationN.png\".\n$d = $word.ActiveDocument\n\nforeach ($n in 1..4) {\n    $oldRef = \"img/Notation$n.png\"\n    $newRef = \"img/Notation/Notation$n.png\"\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldRef\n    $find.Replacement.Text = $newRef\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
